$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.691.49'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '1.656.03'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'1.002"
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').Value = "'302.78"
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').Value = "'51.13"
$ws.Range('E9').Value = '  -1.75%  '
$ws.Range('E10').Value = '  -0.23%  '
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').Value = "'6.450"
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').Value = "'7.444"
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('D16').Value = "'0.00001221"
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('D17').Value = '1.653.95'
$ws.Range('E17').Value = '  +1.00%  '
$ws.Range('D18').Value = "'97.63"
$ws.Range('E18').Value = '  +2.58%  '
$ws.Range('D19').Value = "'0.07040"
$ws.Range('E19').Value = '  +1.00%  '
$ws.Range('D20').Value = "'6.791"
$ws.Range('E20').Value = '  +3.05%  '
$ws.Range('E21').Value = '  +0.54%  '
$ws.Range('D22').Value = "'1.002"
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = "'12.71"
$ws.Range('D24').Value = '23.700.87'
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('D25').Value = "'2.492"
$ws.Range('E25').Value = '  -1.75%  '
$ws.Range('D26').Value = "'3.023"
$ws.Range('E26').Value = '  -1.27%  '
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('D28').Value = "'153.67"
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('D29').Value = "'5.245"
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('D30').Value = "'133.97"
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('D31').Value = '1.839.56'
$ws.Range('E31').Value = '  +1.08%  '
$ws.Range('D32').Value = "'7.113"
$ws.Range('E32').Value = '  +9.51%  '
$ws.Range('D33').Value = "'2.253"
$ws.Range('E33').Value = '  +4.71%  '
$ws.Range('D34').Value = "'12.05"
$ws.Range('E34').Value = '  +4.93%  '
$ws.Range('D35').Value = "'1.055"
$ws.Range('E35').Value = '  -2.78%  '
$ws.Range('D36').Value = "'0.02808"
$ws.Range('E36').Value = '  +1.49%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').Value = "'0.2503"
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').Value = "'0.08811"
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('E39').Value = '  +2.15%  '
$ws.Range('D40').Value = "'0.06986"
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('D41').Value = "'13.03"
$ws.Range('E41').Value = '  +6.51%  '
$ws.Range('D42').Value = "'0.6985"
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').Value = "'15.94"
$ws.Range('E44').Value = '  +2.08%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('D48').Value = "'3.963"
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('D49').Value = "'0.07899"
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('D50').Value = "'128.16"
$ws.Range('E50').Value = '  -0.50%  '
$ws.Range('D51').Value = "'1.177"
$ws.Range('E51').Value = '  -1.18%  '
